# Auto-generated edit script applying odds corrections per the commit diff.
# Sets updated numeric values for the affected cells, grouped by row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: 19 odds updated
$ws.Range("I3").Value = 1.27
$ws.Range("N3").Value = 19
$ws.Range("O3").Value = 1.13
$ws.Range("P3").Value = 6
$ws.Range("Q3").Value = 1.4
$ws.Range("R3").Value = 3
$ws.Range("S3").Value = 1.22
$ws.Range("T3").Value = 4
$ws.Range("U3").Value = 1.75
$ws.Range("V3").Value = 2
$ws.Range("AD3").Value = 11
$ws.Range("AF3").Value = 51
$ws.Range("AH3").Value = 10
$ws.Range("AI3").Value = 8
$ws.Range("AT3").Value = 4
$ws.Range("AV3").Value = 51
$ws.Range("AW3").Value = 3.5
$ws.Range("AY3").Value = 15
$ws.Range("BC3").Value = 451

# Row 4: 13 odds updated
$ws.Range("G4").Value = 1.91
$ws.Range("I4").Value = 4.1
$ws.Range("L4").Value = 4.5
$ws.Range("O4").Value = 1.33
$ws.Range("P4").Value = 3.4
$ws.Range("S4").Value = 1.41
$ws.Range("T4").Value = 2.62
$ws.Range("Y4").Value = 9
$ws.Range("AG4").Value = 301
$ws.Range("AH4").Value = 10
$ws.Range("AJ4").Value = 13
$ws.Range("AP4").Value = 21
$ws.Range("AV4").Value = 51

# Row 6: 23 odds updated
$ws.Range("G6").Value = 2.05
$ws.Range("I6").Value = 3.6
$ws.Range("J6").Value = 2.63
$ws.Range("L6").Value = 4
$ws.Range("N6").Value = 12
$ws.Range("O6").Value = 1.25
$ws.Range("P6").Value = 4
$ws.Range("Q6").Value = 1.8
$ws.Range("R6").Value = 2
$ws.Range("S6").Value = 1.36
$ws.Range("T6").Value = 3
$ws.Range("U6").Value = 1.67
$ws.Range("V6").Value = 2.1
$ws.Range("W6").Value = 8.5
$ws.Range("AC6").Value = 12
$ws.Range("AE6").Value = 13
$ws.Range("AK6").Value = 41
$ws.Range("AL6").Value = 26
$ws.Range("AP6").Value = 21
$ws.Range("AT6").Value = 3
$ws.Range("AX6").Value = 19
$ws.Range("AY6").Value = 26
$ws.Range("BA6").Value = 81

# Row 8: 25 odds updated
$ws.Range("G8").Value = 1.85
$ws.Range("H8").Value = 3.3
$ws.Range("I8").Value = 4.75
$ws.Range("J8").Value = 2.6
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 7
$ws.Range("O8").Value = 1.44
$ws.Range("P8").Value = 2.75
$ws.Range("Q8").Value = 2.38
$ws.Range("R8").Value = 1.57
$ws.Range("S8").Value = 1.53
$ws.Range("T8").Value = 2.38
$ws.Range("U8").Value = 2.1
$ws.Range("V8").Value = 1.67
$ws.Range("Y8").Value = 9
$ws.Range("AB8").Value = 34
$ws.Range("AC8").Value = 7
$ws.Range("AH8").Value = 9.5
$ws.Range("AJ8").Value = 15
$ws.Range("AN8").Value = 3.6
$ws.Range("AP8").Value = 26
$ws.Range("AT8").Value = 2.38
$ws.Range("BB8").Value = 401
$ws.Range("BD8").Value = 126

# Row 9: 26 odds updated
$ws.Range("G9").Value = 2
$ws.Range("I9").Value = 3.9
$ws.Range("J9").Value = 2.75
$ws.Range("L9").Value = 4.75
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 8
$ws.Range("U9").Value = 2.05
$ws.Range("V9").Value = 1.7
$ws.Range("W9").Value = 6
$ws.Range("X9").Value = 8.5
$ws.Range("Z9").Value = 17
$ws.Range("AC9").Value = 7.5
$ws.Range("AH9").Value = 9
$ws.Range("AI9").Value = 19
$ws.Range("AJ9").Value = 13
$ws.Range("AK9").Value = 41
$ws.Range("AL9").Value = 34
$ws.Range("AM9").Value = 41
$ws.Range("AO9").Value = 12
$ws.Range("AS9").Value = 201
$ws.Range("AW9").Value = 5.5
$ws.Range("AX9").Value = 23
$ws.Range("AY9").Value = 34
$ws.Range("AZ9").Value = 81
$ws.Range("BA9").Value = 126
$ws.Range("BB9").Value = 301

# Row 16: 46 odds updated
$ws.Range("G16").Value = 3.55
$ws.Range("H16").Value = 3.15
$ws.Range("I16").Value = 2.07
$ws.Range("J16").Value = 4.1
$ws.Range("K16").Value = 2.05
$ws.Range("L16").Value = 2.65
$ws.Range("M16").Value = 1.09
$ws.Range("N16").Value = 6.2
$ws.Range("O16").Value = 1.4
$ws.Range("P16").Value = 2.72
$ws.Range("Q16").Value = 2.18
$ws.Range("R16").Value = 1.62
$ws.Range("S16").Value = 1.44
$ws.Range("T16").Value = 2.57
$ws.Range("U16").Value = 1.93
$ws.Range("V16").Value = 1.78
$ws.Range("W16").Value = 8.75
$ws.Range("X16").Value = 17.5
$ws.Range("Y16").Value = 12.5
$ws.Range("Z16").Value = 50
$ws.Range("AA16").Value = 37
$ws.Range("AB16").Value = 50
$ws.Range("AC16").Value = 6.2
$ws.Range("AD16").Value = 6.1
$ws.Range("AE16").Value = 16
$ws.Range("AF16").Value = 90
$ws.Range("AG16").Value = 800
$ws.Range("AH16").Value = 6.4
$ws.Range("AI16").Value = 9.25
$ws.Range("AJ16").Value = 8.75
$ws.Range("AK16").Value = 19
$ws.Range("AL16").Value = 18
$ws.Range("AN16").Value = 5.4
$ws.Range("AO16").Value = 20
$ws.Range("AP16").Value = 29
$ws.Range("AQ16").Value = 110
$ws.Range("AR16").Value = 150
$ws.Range("AS16").Value = 400
$ws.Range("AT16").Value = 2.57
$ws.Range("AU16").Value = 7.3
$ws.Range("AV16").Value = 70
$ws.Range("AW16").Value = 3.9
$ws.Range("AX16").Value = 10.5
$ws.Range("AY16").Value = 20
$ws.Range("AZ16").Value = 40
$ws.Range("BA16").Value = 80

# Row 21: 39 odds updated
$ws.Range("G21").Value = 1.75
$ws.Range("H21").Value = 3.75
$ws.Range("I21").Value = 4.2
$ws.Range("J21").Value = 2.38
$ws.Range("K21").Value = 2.3
$ws.Range("L21").Value = 4.5
$ws.Range("M21").Value = 1.04
$ws.Range("N21").Value = 13
$ws.Range("O21").Value = 1.22
$ws.Range("P21").Value = 4
$ws.Range("Q21").Value = 1.7
$ws.Range("R21").Value = 2.1
$ws.Range("U21").Value = 1.67
$ws.Range("V21").Value = 2.1
$ws.Range("W21").Value = 8.5
$ws.Range("X21").Value = 9
$ws.Range("Z21").Value = 15
$ws.Range("AA21").Value = 13
$ws.Range("AD21").Value = 7.5
$ws.Range("AE21").Value = 13
$ws.Range("AF21").Value = 41
$ws.Range("AG21").Value = 151
$ws.Range("AH21").Value = 15
$ws.Range("AI21").Value = 23
$ws.Range("AJ21").Value = 15
$ws.Range("AK21").Value = 41
$ws.Range("AL21").Value = 34
$ws.Range("AM21").Value = 34
$ws.Range("AN21").Value = 4
$ws.Range("AO21").Value = 9
$ws.Range("AP21").Value = 19
$ws.Range("AQ21").Value = 29
$ws.Range("AU21").Value = 7.5
$ws.Range("AW21").Value = 6
$ws.Range("AX21").Value = 21
$ws.Range("AY21").Value = 26
$ws.Range("AZ21").Value = 67
$ws.Range("BA21").Value = 81
$ws.Range("BB21").Value = 151

# Row 25: 11 odds updated
$ws.Range("G25").Value = 2.45
$ws.Range("H25").Value = 3.25
$ws.Range("I25").Value = 2.88
$ws.Range("J25").Value = 3.2
$ws.Range("W25").Value = 8
$ws.Range("X25").Value = 12
$ws.Range("Y25").Value = 10
$ws.Range("AH25").Value = 8.5
$ws.Range("AI25").Value = 13
$ws.Range("AN25").Value = 4.5
$ws.Range("AO25").Value = 15

# Row 50: 9 odds updated
$ws.Range("I50").Value = 9
$ws.Range("L50").Value = 9.5
$ws.Range("O50").Value = 1.17
$ws.Range("P50").Value = 5
$ws.Range("U50").Value = 2.1
$ws.Range("V50").Value = 1.67
$ws.Range("W50").Value = 7.5
$ws.Range("Y50").Value = 9.5
$ws.Range("AY50").Value = 41

